$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 9.072075465861674
$ws.Cells.Item(2, 3).Value = 5.459457105089276
$ws.Cells.Item(2, 4).Value = 5.030686024327976
$ws.Cells.Item(2, 5).Value = 12.60419834982525
$ws.Cells.Item(2, 6).Value = 24.92074885364288
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 11).Value = 8.318395326212755
$ws.Cells.Item(2, 13).Value = 13.55570814820478
$ws.Cells.Item(2, 14).Value = 18.73536573970178
$ws.Cells.Item(2, 15).Value = 22.30085422391772
$ws.Cells.Item(3, 2).Value = 8.793952934428134
$ws.Cells.Item(3, 3).Value = 5.337939148169292
$ws.Cells.Item(3, 4).Value = 4.980291331590565
$ws.Cells.Item(3, 5).Value = 12.38848904863247
$ws.Cells.Item(3, 6).Value = 24.92046979955428
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 11).Value = 8.131818729482911
$ws.Cells.Item(3, 13).Value = 13.38679700773061
$ws.Cells.Item(3, 14).Value = 18.79653686691419
$ws.Cells.Item(3, 15).Value = 22.34993385733485
$ws.Cells.Item(4, 2).Value = 8.620428974443067
$ws.Cells.Item(4, 3).Value = 5.261101922172317
$ws.Cells.Item(4, 4).Value = 4.948568849471407
$ws.Cells.Item(4, 5).Value = 12.25847303576049
$ws.Cells.Item(4, 6).Value = 24.9268184077711
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 11).Value = 8.016568215339728
$ws.Cells.Item(4, 13).Value = 13.28538495281219
$ws.Cells.Item(4, 14).Value = 18.8358296404732
$ws.Cells.Item(4, 15).Value = 22.3848849196188
$ws.Cells.Item(5, 2).Value = 8.549139841634704
$ws.Cells.Item(5, 3).Value = 5.229256138852947
$ws.Cells.Item(5, 4).Value = 4.935451079583808
$ws.Cells.Item(5, 5).Value = 12.20617840493602
$ws.Cells.Item(5, 6).Value = 24.9310449206821
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 11).Value = 7.969497817851945
$ws.Cells.Item(5, 13).Value = 13.24468607267562
$ws.Cells.Item(5, 14).Value = 18.85227894661525
$ws.Cells.Item(5, 15).Value = 22.40033627026939
$ws.Cells.Item(6, 2).Value = 8.537271056031519
$ws.Cells.Item(6, 3).Value = 5.223936691972641
$ws.Cells.Item(6, 4).Value = 4.933261535356757
$ws.Cells.Item(6, 5).Value = 12.19753869556982
$ws.Cells.Item(6, 6).Value = 24.93184572232482
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 11).Value = 7.961677542968665
$ws.Cells.Item(6, 13).Value = 13.23796735225356
$ws.Cells.Item(6, 14).Value = 18.85503678404605
$ws.Cells.Item(6, 15).Value = 22.40297486840617
$ws.Cells.Item(7, 2).Value = 8.619469719911013
$ws.Cells.Item(7, 3).Value = 5.260674566190845
$ws.Cells.Item(7, 4).Value = 4.948392703018856
$ws.Cells.Item(7, 5).Value = 12.25776488587947
$ws.Cells.Item(7, 6).Value = 24.92686877105338
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 11).Value = 8.015933739199028
$ws.Cells.Item(7, 13).Value = 13.28483346955447
$ws.Cells.Item(7, 14).Value = 18.83604970973238
$ws.Cells.Item(7, 15).Value = 22.38508841223688
$ws.Cells.Item(8, 2).Value = 8.976817951665012
$ws.Cells.Item(8, 3).Value = 5.418034233350926
$ws.Cells.Item(8, 4).Value = 5.013475085774296
$ws.Cells.Item(8, 5).Value = 12.52936370980555
$ws.Cells.Item(8, 6).Value = 24.91929967815922
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 11).Value = 8.25424797069582
$ws.Cells.Item(8, 13).Value = 13.49701995160024
$ws.Cells.Item(8, 14).Value = 18.75609854685852
$ws.Cells.Item(8, 15).Value = 22.31677588811809
$ws.Cells.Item(9, 2).Value = 9.651025540820914
$ws.Cells.Item(9, 3).Value = 5.707926457291661
$ws.Cells.Item(9, 4).Value = 5.134677680861452
$ws.Cells.Item(9, 5).Value = 13.07784817311715
$ws.Cells.Item(9, 6).Value = 24.95614923766875
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 11).Value = 8.713189693693664
$ws.Cells.Item(9, 13).Value = 13.9291645895561
$ws.Cells.Item(9, 14).Value = 18.61300869739638
$ws.Cells.Item(9, 15).Value = 22.22113129400507
$ws.Cells.Item(10, 2).Value = 10.12450650574267
$ws.Cells.Item(10, 3).Value = 5.908240153230782
$ws.Cells.Item(10, 4).Value = 5.21950481194796
$ws.Cells.Item(10, 5).Value = 13.48596245894419
$ws.Cells.Item(10, 6).Value = 25.01461590410904
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 11).Value = 9.041585856328867
$ws.Cells.Item(10, 13).Value = 14.25345350989054
$ws.Cells.Item(10, 14).Value = 18.51614322706785
$ws.Cells.Item(10, 15).Value = 22.17434832336486
$ws.Cells.Item(11, 2).Value = 10.33413131093119
$ws.Cells.Item(11, 3).Value = 5.996378294839326
$ws.Cells.Item(11, 4).Value = 5.25711664657133
$ws.Cells.Item(11, 5).Value = 13.67180371275995
$ws.Cells.Item(11, 6).Value = 25.04798163833399
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 11).Value = 9.1883655165755
$ws.Cells.Item(11, 13).Value = 14.40181944346939
$ws.Cells.Item(11, 14).Value = 18.4738524648921
$ws.Cells.Item(11, 15).Value = 22.15818765175485
$ws.Cells.Item(12, 2).Value = 10.41260754905401
$ws.Cells.Item(12, 3).Value = 6.029306675219178
$ws.Cells.Item(12, 4).Value = 5.271213807230005
$ws.Cells.Item(12, 5).Value = 13.74212883153397
$ws.Cells.Item(12, 6).Value = 25.0615836463944
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 11).Value = 9.243518024647152
$ws.Cells.Item(12, 13).Value = 14.45807078135602
$ws.Cells.Item(12, 14).Value = 18.45809177892035
$ws.Cells.Item(12, 15).Value = 22.1528057358506
$ws.Cells.Item(13, 2).Value = 10.39574761277598
$ws.Cells.Item(13, 3).Value = 6.022235142787358
$ws.Cells.Item(13, 4).Value = 5.268184289387937
$ws.Cells.Item(13, 5).Value = 13.72698641358883
$ws.Cells.Item(13, 6).Value = 25.05861130896146
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 11).Value = 9.231659885792523
$ws.Cells.Item(13, 13).Value = 14.44595384402784
$ws.Cells.Item(13, 14).Value = 18.46147485047704
$ws.Cells.Item(13, 15).Value = 22.15393199634905
$ws.Cells.Item(14, 2).Value = 10.34060612129979
$ws.Cells.Item(14, 3).Value = 5.999096402968407
$ws.Cells.Item(14, 4).Value = 5.258279376270731
$ws.Cells.Item(14, 5).Value = 13.67759080802442
$ws.Cells.Item(14, 6).Value = 25.04908134088673
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 11).Value = 9.19291185817038
$ws.Cells.Item(14, 13).Value = 14.40644620549683
$ws.Cells.Item(14, 14).Value = 18.47255074172415
$ws.Cells.Item(14, 15).Value = 22.15773008537411
$ws.Cells.Item(15, 2).Value = 10.30671052323036
$ws.Cells.Item(15, 3).Value = 5.984864445866256
$ws.Cells.Item(15, 4).Value = 5.252193214453197
$ws.Cells.Item(15, 5).Value = 13.64732599594876
$ws.Cells.Item(15, 6).Value = 25.04336970398225
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 11).Value = 9.169120052767193
$ws.Cells.Item(15, 13).Value = 14.3822539798039
$ws.Cells.Item(15, 14).Value = 18.47936807678274
$ws.Cells.Item(15, 15).Value = 22.16015264049479
$ws.Cells.Item(16, 2).Value = 10.11068442850404
$ws.Cells.Item(16, 3).Value = 5.902418442252207
$ws.Cells.Item(16, 4).Value = 5.217026692101729
$ws.Cells.Item(16, 5).Value = 13.47381452064807
$ws.Cells.Item(16, 6).Value = 25.01257107683998
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 11).Value = 9.031936201488847
$ws.Cells.Item(16, 13).Value = 14.2437698896828
$ws.Cells.Item(16, 14).Value = 18.51894262182152
$ws.Cells.Item(16, 15).Value = 22.17550763552608
$ws.Cells.Item(17, 2).Value = 9.988896360876407
$ws.Cells.Item(17, 3).Value = 5.851062311722158
$ws.Cells.Item(17, 4).Value = 5.19519931679616
$ws.Cells.Item(17, 5).Value = 13.36736407085973
$ws.Cells.Item(17, 6).Value = 24.99540673929832
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 11).Value = 8.947069981301093
$ws.Cells.Item(17, 13).Value = 14.15899276365411
$ws.Cells.Item(17, 14).Value = 18.54367380973258
$ws.Cells.Item(17, 15).Value = 22.18624014414004
$ws.Cells.Item(18, 2).Value = 9.918308523074717
$ws.Cells.Item(18, 3).Value = 5.821243604732414
$ws.Cells.Item(18, 4).Value = 5.182553181078505
$ws.Cells.Item(18, 5).Value = 13.3061586651367
$ws.Cells.Item(18, 6).Value = 24.98617183988176
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 11).Value = 8.898014581437454
$ws.Cells.Item(18, 13).Value = 14.11031389705279
$ws.Cells.Item(18, 14).Value = 18.55806556482421
$ws.Cells.Item(18, 15).Value = 22.19289516244565
$ws.Cells.Item(19, 2).Value = 9.894318619669617
$ws.Cells.Item(19, 3).Value = 5.811099991834833
$ws.Cells.Item(19, 4).Value = 5.17825584475158
$ws.Cells.Item(19, 5).Value = 13.28544171903561
$ws.Cells.Item(19, 6).Value = 24.98315473653639
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 11).Value = 8.881365325745367
$ws.Cells.Item(19, 13).Value = 14.09384790722105
$ws.Cells.Item(19, 14).Value = 18.56296708908597
$ws.Cells.Item(19, 15).Value = 22.19523116825374
$ws.Cells.Item(20, 2).Value = 10.00191721818545
$ws.Cells.Item(20, 3).Value = 5.856558377772289
$ws.Cells.Item(20, 4).Value = 5.197532397380261
$ws.Cells.Item(20, 5).Value = 13.37869417740274
$ws.Cells.Item(20, 6).Value = 24.99716796734034
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 11).Value = 8.956129650392917
$ws.Cells.Item(20, 13).Value = 14.16800926310417
$ws.Cells.Item(20, 14).Value = 18.54102385405698
$ws.Cells.Item(20, 15).Value = 22.18504775798571
$ws.Cells.Item(21, 2).Value = 10.35682761232681
$ws.Cells.Item(21, 3).Value = 6.005905098423481
$ws.Cells.Item(21, 4).Value = 5.261192682288145
$ws.Cells.Item(21, 5).Value = 13.69210141368791
$ws.Cells.Item(21, 6).Value = 25.0518543314005
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 11).Value = 9.204305172023512
$ws.Cells.Item(21, 13).Value = 14.41804911964129
$ws.Cells.Item(21, 14).Value = 18.46929060223746
$ws.Cells.Item(21, 15).Value = 22.15659446208823
$ws.Cells.Item(22, 2).Value = 10.58348270643639
$ws.Cells.Item(22, 3).Value = 6.10089529963492
$ws.Cells.Item(22, 4).Value = 5.301947245523525
$ws.Cells.Item(22, 5).Value = 13.89661285547333
$ws.Cells.Item(22, 6).Value = 25.09322848717161
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 11).Value = 9.363974696490027
$ws.Cells.Item(22, 13).Value = 14.58183928465981
$ws.Cells.Item(22, 14).Value = 18.423888210822
$ws.Cells.Item(22, 15).Value = 22.1422993447554
$ws.Cells.Item(23, 2).Value = 10.46302035129782
$ws.Cells.Item(23, 3).Value = 6.050442311631133
$ws.Cells.Item(23, 4).Value = 5.280275324058191
$ws.Cells.Item(23, 5).Value = 13.78751470104574
$ws.Cells.Item(23, 6).Value = 25.07063322034979
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 11).Value = 9.279004340172341
$ws.Cells.Item(23, 13).Value = 14.49440440430087
$ws.Cells.Item(23, 14).Value = 18.44798533578674
$ws.Cells.Item(23, 15).Value = 22.14953502809092
$ws.Cells.Item(24, 2).Value = 9.996032259007585
$ws.Cells.Item(24, 3).Value = 5.854074518758075
$ws.Cells.Item(24, 4).Value = 5.196477914273091
$ws.Cells.Item(24, 5).Value = 13.37357184794311
$ws.Cells.Item(24, 6).Value = 24.99636974358033
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 11).Value = 8.952034594218294
$ws.Cells.Item(24, 13).Value = 14.16393271076893
$ws.Cells.Item(24, 14).Value = 18.54222135857185
$ws.Cells.Item(24, 15).Value = 22.18558532586302
$ws.Cells.Item(25, 2).Value = 9.472096331963352
$ws.Cells.Item(25, 3).Value = 5.631641703196965
$ws.Cells.Item(25, 4).Value = 5.102608798973709
$ws.Cells.Item(25, 5).Value = 12.92824092270033
$ws.Cells.Item(25, 6).Value = 24.94065427313298
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 11).Value = 8.590324488098108
$ws.Cells.Item(25, 13).Value = 13.81084798179818
$ws.Cells.Item(25, 14).Value = 18.6502611353412
$ws.Cells.Item(25, 15).Value = 22.24288933889773
